$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1675.5
$ws.Cells.Item(12, 9).Value = 1900.3334
$ws.Cells.Item(12, 11).Value = 1900.3334
$ws.Cells.Item(12, 13).Value = -1730.3334
$ws.Cells.Item(33, 8).Value = 362.9091
$ws.Cells.Item(33, 9).Value = 297
$ws.Cells.Item(33, 10).Value = 732
$ws.Cells.Item(33, 11).Value = 297
$ws.Cells.Item(33, 12).Value = 732
$ws.Cells.Item(33, 13).Value = -68
$ws.Cells.Item(33, 14).Value = -1190
$ws.Cells.Item(38, 8).Value = 495.2
$ws.Cells.Item(38, 9).Value = 242
$ws.Cells.Item(38, 11).Value = 726
$ws.Cells.Item(38, 13).Value = -354
$ws.Cells.Item(39, 8).Value = 77.42856999999999
$ws.Cells.Item(39, 9).Value = 77.42856999999999
$ws.Cells.Item(39, 11).Value = 232.28571
$ws.Cells.Item(39, 13).Value = 63.71429000000001
$ws.Cells.Item(97, 8).Value = 1082
$ws.Cells.Item(97, 10).Value = 1082
$ws.Cells.Item(97, 12).Value = 3246
$ws.Cells.Item(97, 14).Value = -4238
$ws.Cells.Item(98, 8).Value = 888.4545000000001
$ws.Cells.Item(98, 9).Value = 829.1429000000001
$ws.Cells.Item(98, 11).Value = 829.1429000000001
$ws.Cells.Item(98, 13).Value = 668.8570999999999
$ws.Cells.Item(111, 8).Value = 4634.1816
$ws.Cells.Item(111, 9).Value = 1546
$ws.Cells.Item(111, 10).Value = 8340
$ws.Cells.Item(111, 11).Value = 4638
$ws.Cells.Item(111, 12).Value = 25020
$ws.Cells.Item(111, 13).Value = -1571
$ws.Cells.Item(111, 14).Value = -31154
$ws.Cells.Item(116, 8).Value = 5496.364
$ws.Cells.Item(116, 9).Value = 1651.6666
$ws.Cells.Item(116, 10).Value = 6938.125
$ws.Cells.Item(116, 11).Value = 1651.6666
$ws.Cells.Item(116, 12).Value = 6938.125
$ws.Cells.Item(116, 13).Value = 1790.3334
$ws.Cells.Item(116, 14).Value = -13822.125
$ws.Cells.Item(122, 8).Value = 888.4545000000001
$ws.Cells.Item(122, 9).Value = 829.1429000000001
$ws.Cells.Item(122, 11).Value = 2487.4287
$ws.Cells.Item(122, 13).Value = -37.42870000000039
$ws.Cells.Item(129, 8).Value = 839.0706
$ws.Cells.Item(129, 9).Value = 796.6667
$ws.Cells.Item(129, 10).Value = 840.62195
$ws.Cells.Item(129, 11).Value = 2390.0001
$ws.Cells.Item(129, 12).Value = 2521.86585
$ws.Cells.Item(129, 13).Value = 2609.9999
$ws.Cells.Item(129, 14).Value = -12521.86585
$ws.Cells.Item(138, 8).Value = 1566.1794
$ws.Cells.Item(138, 10).Value = 3589.923
$ws.Cells.Item(138, 12).Value = 10769.769
$ws.Cells.Item(138, 14).Value = -21049.769
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20102.727
$ws.Cells.Item(32, 9).Value = 21795.32
$ws.Cells.Item(32, 10).Value = 3176.8
$ws.Cells.Item(32, 11).Value = 21795.32
$ws.Cells.Item(32, 12).Value = 3176.8
$ws.Cells.Item(32, 13).Value = -21508.32
$ws.Cells.Item(32, 14).Value = -3750.8
$ws.Cells.Item(102, 8).Value = 2566
$ws.Cells.Item(102, 9).Value = 1031.3846
$ws.Cells.Item(102, 10).Value = 6556
$ws.Cells.Item(102, 11).Value = 1031.3846
$ws.Cells.Item(102, 12).Value = 6556
$ws.Cells.Item(102, 13).Value = 590.6153999999999
$ws.Cells.Item(102, 14).Value = -9800
$ws.Cells.Item(110, 8).Value = 4986.6665
$ws.Cells.Item(110, 9).Value = 4980
$ws.Cells.Item(110, 10).Value = 5000
$ws.Cells.Item(110, 11).Value = 4980
$ws.Cells.Item(110, 12).Value = 5000
$ws.Cells.Item(110, 13).Value = -2935
$ws.Cells.Item(110, 14).Value = -9090
$ws.Cells.Item(122, 8).Value = 2251.8333
$ws.Cells.Item(122, 9).Value = 2251.8333
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 6755.499899999999
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -4305.499899999999
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 18616.967
$ws.Cells.Item(132, 9).Value = 1966.5
$ws.Cells.Item(132, 11).Value = 5899.5
$ws.Cells.Item(132, 13).Value = -3369.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1102.28
$ws.Cells.Item(20, 9).Value = 1168.1765
$ws.Cells.Item(20, 10).Value = 962.25
$ws.Cells.Item(20, 11).Value = 1168.1765
$ws.Cells.Item(20, 12).Value = 962.25
$ws.Cells.Item(20, 13).Value = -921.1765
$ws.Cells.Item(20, 14).Value = -1456.25
$ws.Cells.Item(86, 8).Value = 2064.4
$ws.Cells.Item(86, 9).Value = 1859.2
$ws.Cells.Item(86, 10).Value = 2680
$ws.Cells.Item(86, 11).Value = 1859.2
$ws.Cells.Item(86, 12).Value = 2680
$ws.Cells.Item(86, 13).Value = -736.2
$ws.Cells.Item(86, 14).Value = -4926
$ws.Cells.Item(89, 8).Value = 2064.4
$ws.Cells.Item(89, 9).Value = 1859.2
$ws.Cells.Item(89, 10).Value = 2680
$ws.Cells.Item(89, 11).Value = 9296
$ws.Cells.Item(89, 12).Value = 13400
$ws.Cells.Item(89, 13).Value = -3680
$ws.Cells.Item(89, 14).Value = -24632
$ws.Cells.Item(103, 8).Value = 20657
$ws.Cells.Item(103, 10).Value = 20657
$ws.Cells.Item(103, 12).Value = 20657
$ws.Cells.Item(103, 14).Value = -23001
$ws.Cells.Item(105, 8).Value = 5002590
$ws.Cells.Item(105, 9).Value = 2699.75
$ws.Cells.Item(105, 10).Value = 8335850
$ws.Cells.Item(105, 11).Value = 2699.75
$ws.Cells.Item(105, 12).Value = 8335850
$ws.Cells.Item(105, 13).Value = -952.75
$ws.Cells.Item(105, 14).Value = -8339344
$ws.Cells.Item(107, 8).Value = 2655.5
$ws.Cells.Item(107, 9).Value = 2655.5
$ws.Cells.Item(107, 11).Value = 2655.5
$ws.Cells.Item(107, 13).Value = -735.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 2471.6667
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 2471.6667
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 2471.6667
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = -2695.6667
$ws.Cells.Item(31, 8).Value = 2446.4517
$ws.Cells.Item(31, 9).Value = 2041.2916
$ws.Cells.Item(31, 10).Value = 3835.5715
$ws.Cells.Item(31, 11).Value = 2041.2916
$ws.Cells.Item(31, 12).Value = 3835.5715
$ws.Cells.Item(31, 13).Value = -1746.2916
$ws.Cells.Item(31, 14).Value = -4425.5715
$ws.Cells.Item(34, 8).Value = 2446.4517
$ws.Cells.Item(34, 9).Value = 2041.2916
$ws.Cells.Item(34, 10).Value = 3835.5715
$ws.Cells.Item(34, 11).Value = 2041.2916
$ws.Cells.Item(34, 12).Value = 3835.5715
$ws.Cells.Item(34, 13).Value = -1839.2916
$ws.Cells.Item(34, 14).Value = -4239.5715
$ws.Cells.Item(58, 8).Value = 27804.895
$ws.Cells.Item(58, 9).Value = 2260.4285
$ws.Cells.Item(58, 11).Value = 2260.4285
$ws.Cells.Item(58, 13).Value = -2057.4285
$ws.Cells.Item(62, 8).Value = 5716.7144
$ws.Cells.Item(62, 9).Value = 4001.6667
$ws.Cells.Item(62, 11).Value = 4001.6667
$ws.Cells.Item(62, 13).Value = -3377.6667
$ws.Cells.Item(65, 8).Value = 5716.7144
$ws.Cells.Item(65, 9).Value = 4001.6667
$ws.Cells.Item(65, 11).Value = 20008.3335
$ws.Cells.Item(65, 13).Value = -16888.3335
$ws.Cells.Item(86, 8).Value = 20644.75
$ws.Cells.Item(86, 9).Value = 3833.3333
$ws.Cells.Item(86, 10).Value = 30731.6
$ws.Cells.Item(86, 11).Value = 3833.3333
$ws.Cells.Item(86, 12).Value = 30731.6
$ws.Cells.Item(86, 13).Value = -2710.3333
$ws.Cells.Item(86, 14).Value = -32977.6
$ws.Cells.Item(89, 8).Value = 20644.75
$ws.Cells.Item(89, 9).Value = 3833.3333
$ws.Cells.Item(89, 10).Value = 30731.6
$ws.Cells.Item(89, 11).Value = 19166.6665
$ws.Cells.Item(89, 12).Value = 153658
$ws.Cells.Item(89, 13).Value = -13550.6665
$ws.Cells.Item(89, 14).Value = -164890
$ws.Cells.Item(94, 8).Value = 2925.2104
$ws.Cells.Item(94, 9).Value = 2219
$ws.Cells.Item(94, 10).Value = 3337.1667
$ws.Cells.Item(94, 11).Value = 2219
$ws.Cells.Item(94, 12).Value = 3337.1667
$ws.Cells.Item(94, 13).Value = -1768
$ws.Cells.Item(94, 14).Value = -4239.1667
$ws.Cells.Item(136, 8).Value = 27804.895
$ws.Cells.Item(136, 9).Value = 2260.4285
$ws.Cells.Item(136, 11).Value = 6781.2855
$ws.Cells.Item(136, 13).Value = -4231.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 112.625
$ws.Cells.Item(4, 9).Value = 114.42857
$ws.Cells.Item(4, 11).Value = 343.28571
$ws.Cells.Item(4, 13).Value = -231.28571
$ws.Cells.Item(17, 8).Value = 499.75
$ws.Cells.Item(17, 9).Value = 166.33333
$ws.Cells.Item(17, 10).Value = 1500
$ws.Cells.Item(17, 11).Value = 498.99999
$ws.Cells.Item(17, 12).Value = 4500
$ws.Cells.Item(17, 13).Value = -329.99999
$ws.Cells.Item(17, 14).Value = -4838
$ws.Cells.Item(34, 8).Value = 786
$ws.Cells.Item(34, 9).Value = 346.66666
$ws.Cells.Item(34, 10).Value = 895.8333
$ws.Cells.Item(34, 11).Value = 1039.99998
$ws.Cells.Item(34, 12).Value = 2687.4999
$ws.Cells.Item(34, 13).Value = -955.9999800000001
$ws.Cells.Item(34, 14).Value = -2855.4999
$ws.Cells.Item(39, 8).Value = 2414.1
$ws.Cells.Item(39, 10).Value = 2471.2222
$ws.Cells.Item(39, 12).Value = 7413.6666
$ws.Cells.Item(39, 14).Value = -8001.6666
$ws.Cells.Item(55, 8).Value = 2720.1667
$ws.Cells.Item(55, 10).Value = 2720.1667
$ws.Cells.Item(55, 12).Value = 8160.500100000001
$ws.Cells.Item(55, 14).Value = -8514.500100000001
$ws.Cells.Item(131, 8).Value = 770.54
$ws.Cells.Item(131, 10).Value = 782.6185
$ws.Cells.Item(131, 12).Value = 2347.8555
$ws.Cells.Item(131, 14).Value = -12427.8555
$ws.Cells.Item(141, 8).Value = 2601.9048
$ws.Cells.Item(141, 9).Value = 2387.6924
$ws.Cells.Item(141, 10).Value = 2950
$ws.Cells.Item(141, 11).Value = 7163.0772
$ws.Cells.Item(141, 12).Value = 8850
$ws.Cells.Item(141, 13).Value = -1983.0772
$ws.Cells.Item(141, 14).Value = -19210
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 9000
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 9000
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 9000
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).Value = -9224
$ws.Cells.Item(97, 8).Value = 2072.04
$ws.Cells.Item(97, 9).Value = 1484.7778
$ws.Cells.Item(97, 10).Value = 3582.1428
$ws.Cells.Item(97, 11).Value = 1484.7778
$ws.Cells.Item(97, 12).Value = 3582.1428
$ws.Cells.Item(97, 13).Value = -988.7778000000001
$ws.Cells.Item(97, 14).Value = -4574.1428
$ws.Cells.Item(132, 8).Value = 118971.16
$ws.Cells.Item(132, 9).Value = 114736.336
$ws.Cells.Item(132, 10).Value = 128499.5
$ws.Cells.Item(132, 11).Value = 344209.008
$ws.Cells.Item(132, 12).Value = 385498.5
$ws.Cells.Item(132, 13).Value = -341679.008
$ws.Cells.Item(132, 14).Value = -390558.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 1133333.4
$ws.Cells.Item(2, 9).Value = 1178571.4
$ws.Cells.Item(2, 11).Value = 1178571.4
$ws.Cells.Item(2, 13).Value = -1178459.4
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 1035260.4
$ws.Cells.Item(122, 9).Value = 2181194
$ws.Cells.Item(122, 10).Value = 3920
$ws.Cells.Item(122, 11).Value = 6543582
$ws.Cells.Item(122, 12).Value = 11760
$ws.Cells.Item(122, 13).Value = -6541132
$ws.Cells.Item(122, 14).Value = -16660
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 46668.668
$ws.Cells.Item(2, 9).Value = 10000
$ws.Cells.Item(2, 11).Value = 10000
$ws.Cells.Item(2, 13).Value = -9888
$ws.Cells.Item(81, 8).Value = 1581
$ws.Cells.Item(81, 9).Value = 1645.5555
$ws.Cells.Item(81, 10).Value = 1000
$ws.Cells.Item(81, 11).Value = 3291.111
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 13).Value = -2230.111
$ws.Cells.Item(81, 14).Value = -4122
$ws.Cells.Item(84, 8).Value = 1581
$ws.Cells.Item(84, 9).Value = 1645.5555
$ws.Cells.Item(84, 10).Value = 1000
$ws.Cells.Item(84, 11).Value = 16455.555
$ws.Cells.Item(84, 12).Value = 10000
$ws.Cells.Item(84, 13).Value = -11151.555
$ws.Cells.Item(84, 14).Value = -20608
